# Edit described by the commit:
#   Slide with SlideID=286 ("느낀 점"), placeholder shape id=9
#   ("텍스트 개체 틀 8"): the first bullet's opening clause
#   "결측 값을 0으로 처리하지 않고" is reworded to "정확도를 높일",
#   while the remainder of the sentence (" 다른 방법을 찾으려 했지만
#   오류로 어려움을 겪었으며, 현재의 방법을 사용하게 된 점이
#   아쉬웠습니다.") is left untouched.

$p = $ppt.ActivePresentation

# Locate the target slide by its stable SlideID (286) rather than a
# positional index, so the script is resilient to slide reordering.
$targetSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    if ($sl.SlideID -eq 286) {
        $targetSlide = $sl
    }
}

# Locate the target shape on that slide by its stable shape id (9)
# rather than a positional index.
$targetShape = $null
for ($i = 1; $i -le $targetSlide.Shapes.Count; $i++) {
    $sh = $targetSlide.Shapes.Item($i)
    if ($sh.Id -eq 9) {
        $targetShape = $sh
    }
}

$tr = $targetShape.TextFrame.TextRange

# Replace the old opening clause "결측 값을 0으로 처리하지 않고"
# (the first 17 characters of the text body) with the new wording
# "정확도를 높일". The rest of the sentence (" 다른 방법을 찾으려
# 했지만 오류로 어려움을 겪었으며, 현재의 방법을 사용하게 된 점이
# 아쉬웠습니다.") is left untouched, so it keeps its original
# run/character formatting (bold, 00B0F0 fill, Söhne font, etc.).
$oldClause = $tr.Characters(1, 17)
$oldClause.Text = "정확도를 높일"
